# Atualização de bases das ligas, do dia: 28-06-2024 às 19:47
#
# Rows 298-304 (match id + odds data, columns B:AD) get cyclically rotated
# "up" by one row: row 298 takes the data that used to be in row 299,
# row 299 takes what used to be in row 300, ... , row 303 takes what used
# to be in row 304, and row 304 takes what used to be (originally) in row 298.
# Column A (the sequential index) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 298
$lastRow  = 304
$startCol = 2   # column B
$endCol   = 30  # column AD

# Snapshot the current values for columns B:AD across the affected rows
# before we start overwriting anything. (Note: reading a COM property in
# this environment requires the explicit getter call syntax "Value()".)
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rng = $ws.Range($ws.Cells.Item($r, $startCol), $ws.Cells.Item($r, $endCol))
    $snapshot[$r] = $rng.Value()
}

# Write each row's new values: row r gets the snapshot from row r+1,
# wrapping the last row back around to the first row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $r + 1
    if ($srcRow -gt $lastRow) {
        $srcRow = $firstRow
    }
    $destRng = $ws.Range($ws.Cells.Item($r, $startCol), $ws.Cells.Item($r, $endCol))
    $destRng.Value = $snapshot[$srcRow]
}
